# Week 6 and 8 evidence additions to the Constraints Plan.
# Appends extra business-impact commentary to the "Effect of Constraint on
# Product/Project" and "Solution" columns for each constraint row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Hardware & Software Platforms
$ws.Range("D5").Value = "Lacks usability identified during design phase as important, from user needs. Cannot be used across different platfrom types, due to not being usable on small screens like phones and tablets. This will have an impact on the potential market and reduce take-up."
$ws.Range("G5").Value = "Needs further development to ensure all types of users are catered for. Different media types can be easily catered for using Flexbox & media queries in CSS."

# Row 6 - Performance Requirements
$ws.Range("D6").Value = "Does not provide enough options for the user to effectively manage their budget. Users will prefer to use competitor product which includes this functionality."
$ws.Range("G6").Value = "Database neds adjusting to include a budget table that could then be manipulated and sorted by a model."

# Row 7 - Persistenet storage & transactions
$ws.Range("D7").Value = "No transactions, but storage may be an issue in the future depending on the volume of users. Not considered to be a major significant business risk at this point, but will become an issue when userbase increases."
$ws.Range("G7").Value = "Database would have to be expanded to accommodate. This could be achieved by changing the Transactions table from INT2 to INT8."

# Row 8 - Usability
$ws.Range("D8").Value = "Uasability can impacted by lack of portability.`nNot accessible to all user types, which will reduce appeal, and have an impact on the volume of users."
$ws.Range("G8").Value = "Website is designed on a tabular layout which has ben made accessible for screen readers. Navigation is provided by large accessible buttons, and mouse clicks are reduced to a minimum. `nFurther work required on portability aspects, as per Hardware & Software constraint above."

# Row 10 - Time
$ws.Range("D10").Value = "4.5 days development time for MVP and extensions was allocated by CodeClan. Schedule very tight. `nBasic product not visually appealing which will affect marketability."

# Update the saved selection to match the author's final cursor position.
$ws.Range("J8").Select()
